$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextCell $ws.Range('D2') '51.067.64'
Set-TextCell $ws.Range('E2') '  +0.36%  '
Set-TextCell $ws.Range('D3') '2.954.95'
Set-TextCell $ws.Range('E3') '  +0.54%  '
Set-TextCell $ws.Range('E4') '  +0.09%  '
Set-TextCell $ws.Range('D5') '379.77'
Set-TextCell $ws.Range('E5') '  +1.20%  '
Set-TextCell $ws.Range('D6') '101.94'
Set-TextCell $ws.Range('E6') '  +0.44%  '
Set-TextCell $ws.Range('D7') '0.543'
Set-TextCell $ws.Range('E8') '  -0.03%  '
Set-TextCell $ws.Range('D9') '0.589'
Set-TextCell $ws.Range('E9') '  +1.07%  '
Set-TextCell $ws.Range('D10') '36.33'
Set-TextCell $ws.Range('E10') '  +0.45%  '
Set-TextCell $ws.Range('E11') '  -1.10%  '
Set-TextCell $ws.Range('D12') '0.0853'
Set-TextCell $ws.Range('E12') '  +1.89%  '
Set-TextCell $ws.Range('D13') '3.418.94'
Set-TextCell $ws.Range('E13') '  +0.41%  '
Set-TextCell $ws.Range('E14') '  +5.13%  '
Set-TextCell $ws.Range('D15') '18.27'
Set-TextCell $ws.Range('E15') '  +2.26%  '
Set-TextCell $ws.Range('B16') 'WrappedEther'
Set-TextCell $ws.Range('C16') 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextCell $ws.Range('D16') '2.948.03'
Set-TextCell $ws.Range('E16') '  +0.69%  '
Set-TextCell $ws.Range('B17') 'Uniswap'
Set-TextCell $ws.Range('C17') 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextCell $ws.Range('D17') '11.18'
Set-TextCell $ws.Range('E17') '  +6.51%  '
Set-TextCell $ws.Range('D18') '0.994'
Set-TextCell $ws.Range('E18') '  +1.69%  '
Set-TextCell $ws.Range('D19') '51.148.93'
Set-TextCell $ws.Range('E19') '  +0.62%  '
Set-TextCell $ws.Range('D20') '3.12'
Set-TextCell $ws.Range('E20') '  -0.38%  '
Set-TextCell $ws.Range('D21') '12.32'
Set-TextCell $ws.Range('D22') '0.0₃0958'
Set-TextCell $ws.Range('E22') '  +0.68%  '
Set-TextCell $ws.Range('B23') 'Litecoin'
Set-TextCell $ws.Range('C23') 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextCell $ws.Range('D23') '70.34'
Set-TextCell $ws.Range('E23') '  +2.91%  '
Set-TextCell $ws.Range('B24') 'PancakeSwap'
Set-TextCell $ws.Range('C24') 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextCell $ws.Range('D24') '3.30'
Set-TextCell $ws.Range('E24') '  +7.52%  '
Set-TextCell $ws.Range('D25') '266.92'
Set-TextCell $ws.Range('E25') '  +1.24%  '
Set-TextCell $ws.Range('D26') '7.81'
Set-TextCell $ws.Range('E26') '  -4.93%  '
Set-TextCell $ws.Range('B27') 'Dai'
Set-TextCell $ws.Range('C27') 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextCell $ws.Range('D27') '0.999'
Set-TextCell $ws.Range('E27') '  -0.07%  '
Set-TextCell $ws.Range('B28') 'RenderToken'
Set-TextCell $ws.Range('C28') 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell $ws.Range('D28') '7.15'
Set-TextCell $ws.Range('E28') '  -9.12%  '
Set-TextCell $ws.Range('B29') 'EthereumClassic'
Set-TextCell $ws.Range('C29') 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextCell $ws.Range('D29') '25.81'
Set-TextCell $ws.Range('E29') '  +0.95%  '
Set-TextCell $ws.Range('B30') 'Kaspa'
Set-TextCell $ws.Range('C30') 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextCell $ws.Range('D30') '0.165'
Set-TextCell $ws.Range('E30') '  -2.56%  '
Set-TextCell $ws.Range('E31') '  +0.21%  '
Set-TextCell $ws.Range('D32') '10.25'
Set-TextCell $ws.Range('E32') '  +3.97%  '
Set-TextCell $ws.Range('B33') 'OKB'
Set-TextCell $ws.Range('C33') 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextCell $ws.Range('D33') '51.01'
Set-TextCell $ws.Range('E33') '  +0.69%  '
Set-TextCell $ws.Range('B34') 'InjectiveProtocol'
Set-TextCell $ws.Range('C34') 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextCell $ws.Range('D34') '34.20'
Set-TextCell $ws.Range('E34') '  +2.14%  '
Set-TextCell $ws.Range('E35') '  +2.05%  '
Set-TextCell $ws.Range('E36') '  -1.42%  '
Set-TextCell $ws.Range('E37') '  +0.00%  '
Set-TextCell $ws.Range('D38') '3.21'
Set-TextCell $ws.Range('E38') '  +6.26%  '
Set-TextCell $ws.Range('E39') '  +0.81%  '
Set-TextCell $ws.Range('D40') '1.82'
Set-TextCell $ws.Range('D41') '16.43'
Set-TextCell $ws.Range('E41') '  +0.79%  '
Set-TextCell $ws.Range('D42') '2.49'
Set-TextCell $ws.Range('E42') '  -0.95%  '
Set-TextCell $ws.Range('D43') '124.66'
Set-TextCell $ws.Range('E43') '  +3.68%  '
Set-TextCell $ws.Range('B44') 'EnergySwap'
Set-TextCell $ws.Range('C44') 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell $ws.Range('D44') '21.41'
Set-TextCell $ws.Range('E44') '  +1.92%  '
Set-TextCell $ws.Range('B45') 'NEARProtocol'
Set-TextCell $ws.Range('C45') 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell $ws.Range('D45') '3.51'
Set-TextCell $ws.Range('E45') '  +7.12%  '
Set-TextCell $ws.Range('D46') '0.271'
Set-TextCell $ws.Range('E46') '  -5.48%  '
Set-TextCell $ws.Range('E47') '  -0.12%  '
Set-TextCell $ws.Range('E48') '  +2.79%  '
Set-TextCell $ws.Range('D49') '2.041.61'
Set-TextCell $ws.Range('E49') '  +2.54%  '
Set-TextCell $ws.Range('E50') '  -5.62%  '
Set-TextCell $ws.Range('E51') '  +6.99%  '
